$wb = $excel.ActiveWorkbook

$wsUsers    = $wb.Worksheets.Item(1)  # "Users"
$wsSummary  = $wb.Worksheets.Item(2)  # "SummaryLogs"
$wsDetail   = $wb.Worksheets.Item(3)  # "DetailLogs"

# ---------------------------------------------------------------------------
# SummaryLogs sheet: the engagement names in column A are being replaced with
# fuller legal-entity names. The previous (shorter) engagement names are kept,
# but moved into a new column E alongside the row they used to head. The old
# 4th data row (Debevoise) is removed entirely from this sheet.
# ---------------------------------------------------------------------------

$oldSummaryA2 = $wsSummary.Range("A2").Value2
$oldSummaryA3 = $wsSummary.Range("A3").Value2

# Move old project names into column E, preserving their original look
# (row2's old label was bold/default-colored like the header; row3's old
# label was bold with explicit black text).
$wsSummary.Range("E2").ClearFormats()
$wsSummary.Range("E2").Value = $oldSummaryA2
$wsSummary.Range("E2").Font.Bold = $true

$wsSummary.Range("E3").ClearFormats()
$wsSummary.Range("E3").Value = $oldSummaryA3
$wsSummary.Range("E3").Font.Bold = $true
$wsSummary.Range("E3").Font.Color = 0

# Replace column A with the new, fuller project names (plain formatting).
$wsSummary.Range("A2").ClearFormats()
$wsSummary.Range("A2").Value = "Eventide-Eventide Asset Management, LLC-FVA-111771"

$wsSummary.Range("A3").ClearFormats()
$wsSummary.Range("A3").Value = "(MEH) Thompson_CLP Toxicology-Thompson Hine, LLP-FVA-26495"

# Remove the old 4th row (Debevoise engagement) from SummaryLogs entirely.
[void]$wsSummary.Rows.Item(4).Delete()

# ---------------------------------------------------------------------------
# DetailLogs sheet: same column-A replacement, but here a new column D is
# added to retain the old project names (all 3 rows keep their old value in
# column D), and a brand new 4th row is introduced for Debevoise with its own
# fuller legal-entity name.
# ---------------------------------------------------------------------------

$oldDetailA2 = $wsDetail.Range("A2").Value2
$oldDetailA3 = $wsDetail.Range("A3").Value2
$oldDetailA4 = $wsDetail.Range("A4").Value2

$wsDetail.Range("D2").Value = $oldDetailA2

$wsDetail.Range("D3").ClearFormats()
$wsDetail.Range("D3").Value = $oldDetailA3
$wsDetail.Range("D3").Font.Bold = $false
$wsDetail.Range("D3").Font.Color = 0

$wsDetail.Range("D4").Value = $oldDetailA4

$wsDetail.Range("A2").ClearFormats()
$wsDetail.Range("A2").Value = "Eventide-Eventide Asset Management, LLC-FVA-111771"

$wsDetail.Range("A3").ClearFormats()
$wsDetail.Range("A3").Value = "(MEH) Thompson_CLP Toxicology-Thompson Hine, LLP-FVA-26495"

$wsDetail.Range("A4").Value = "Debevoise_Xie (Consulting)-Debevoise & Plimpton LLP-FVA-26378"

# ---------------------------------------------------------------------------
# Selection / active-sheet state, matching the recorded UI state of the edit:
# Users tab is no longer the active tab; SummaryLogs becomes active with the
# whole 4th row selected (about to be deleted/was just worked on); DetailLogs
# keeps a parked selection at D7.
# ---------------------------------------------------------------------------

[void]$wsUsers.Range("A2").Select()
[void]$wsDetail.Range("D7").Select()

[void]$wsSummary.Activate()
[void]$wsSummary.Range("A4:XFD4").Select()
